$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RECURSO" (D) and "ERA" (E) columns entirely.
$ws.Range("D1:E1").EntireColumn.Delete()

# Match the author's resulting selection (was E4, now C4 after the
# column removal shifted the active cell).
[void]$ws.Range("C4").Select()
